$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- "Metadata" sheet: bump the "Date" value ---
$ws1.Range("B8").Value = "2024-03-19T14:09:21+00:00"

# --- "Include from CGH abnormalitie" sheet: turn the old one-line
#     "Codes" / "All codes" summary into a real Concept table.
#
#     Before:            After:
#       Codes               Concept | Description
#       All codes           CNV0    | CNVs not explaining phenotypes
#       (blank) (blank)     CNV1    | CNVs explaining phenotypes
#       System URI | <url>  (blank) | (blank)
#                           System URI | <url>
#
#     Insert a fresh row above the existing blank separator row so the
#     separator row and the "System URI" row below it simply slide down
#     one position, keeping their original cell content/format intact.
$ws2.Rows(3).Insert()

# Row 1 gets a second column (header); copy formatting from A1 (header style).
$ws2.Range("A1").Copy() | Out-Null
$ws2.Range("B1").PasteSpecial(-4122) | Out-Null

# Row 2 gets a second column; copy formatting from A2 (body style).
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("B2").PasteSpecial(-4122) | Out-Null

# New row 3 becomes the "CNV1" row; copy formatting from row 2 (now fully
# styled across both columns).
$ws2.Range("A2:B2").Copy() | Out-Null
$ws2.Range("A3:B3").PasteSpecial(-4122) | Out-Null

# Fill in the new values.
$ws2.Range("A1").Value = "Concept"
$ws2.Range("B1").Value = "Description"

$ws2.Range("A2").Value = "CNV0"
$ws2.Range("B2").Value = "CNVs not explaining phenotypes"

$ws2.Range("A3").Value = "CNV1"
$ws2.Range("B3").Value = "CNVs explaining phenotypes"
